$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text (daily conversion rates) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.33 = 8796.84 pesos`n✅ 8796.84 pesos = 2.33 = 969.09 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update tasas sheet rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 429.7
$ws2.Range("O10").Value = 3780
$ws2.Range("N12").Value = 3780
$ws2.Range("O12").Value = 416.42
